$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Source cells used to copy exact cell formatting (styles) for columns
# B (Composition), C (Structure), D (Processing when "AAM"), E (Material Comment)
$srcB = $ws.Range("B288")
$srcC = $ws.Range("C81")
$srcD48 = $ws.Range("D165")
$srcE = $ws.Range("E288")

# Row 305
$ws.Range("A305").Value = "Al5"
$srcB.Copy($ws.Range("B305"))
$ws.Range("B305").Value = "Al5 Sc25 Ti35 Zr35"
$srcC.Copy($ws.Range("C305"))
$ws.Range("C305").Value = "HCP"
$srcD48.Copy($ws.Range("D305"))
$ws.Range("D305").Value = "AAM"
$ws.Range("F305").Value = "compressive yield stress"
$ws.Range("G305").Value = "EXP"
$ws.Range("I305").Value = 298
$ws.Range("J305").Value = 1615000000
$ws.Range("K305").Value = 21000000
$ws.Range("L305").Value = "Pa"
$ws.Range("M305").Value = "T3"
$ws.Range("N305").Value = "10.1016/j.matchar.2024.113730"

# Row 306
$ws.Range("A306").Value = "Al10"
$srcB.Copy($ws.Range("B306"))
$ws.Range("B306").Value = "Al10 Sc20 Ti35 Zr35"
$srcC.Copy($ws.Range("C306"))
$ws.Range("C306").Value = "HCP+BCC"
$srcD48.Copy($ws.Range("D306"))
$ws.Range("D306").Value = "AAM"
$ws.Range("F306").Value = "compressive yield stress"
$ws.Range("G306").Value = "EXP"
$ws.Range("I306").Value = 298
$ws.Range("J306").Value = 1268000000
$ws.Range("K306").Value = 3000000
$ws.Range("L306").Value = "Pa"
$ws.Range("M306").Value = "T3"
$ws.Range("N306").Value = "10.1016/j.matchar.2024.113730"

# Row 307
$ws.Range("A307").Value = "Al15"
$srcB.Copy($ws.Range("B307"))
$ws.Range("B307").Value = "Al15 Sc15 Ti35 Zr35"
$srcC.Copy($ws.Range("C307"))
$ws.Range("C307").Value = "BCC"
$srcD48.Copy($ws.Range("D307"))
$ws.Range("D307").Value = "AAM"
$ws.Range("F307").Value = "compressive yield stress"
$ws.Range("G307").Value = "EXP"
$ws.Range("I307").Value = 298
$ws.Range("J307").Value = 2007000000
$ws.Range("K307").Value = 59000000
$ws.Range("L307").Value = "Pa"
$ws.Range("M307").Value = "T3"
$ws.Range("N307").Value = "10.1016/j.matchar.2024.113730"

# Row 308
$ws.Range("A308").Value = "Al5"
$srcB.Copy($ws.Range("B308"))
$ws.Range("B308").Value = "Al5 Sc25 Ti35 Zr35"
$srcC.Copy($ws.Range("C308"))
$ws.Range("C308").Value = "HCP+HCP"
$ws.Range("D308").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E308"))
$ws.Range("E308").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F308").Value = "compressive yield stress"
$ws.Range("G308").Value = "EXP"
$ws.Range("I308").Value = 298
$ws.Range("J308").Value = 1351000000
$ws.Range("K308").Value = 24000000
$ws.Range("L308").Value = "Pa"
$ws.Range("M308").Value = "T3"
$ws.Range("N308").Value = "10.1016/j.matchar.2024.113730"

# Row 309
$ws.Range("A309").Value = "Al10"
$srcB.Copy($ws.Range("B309"))
$ws.Range("B309").Value = "Al10 Sc20 Ti35 Zr35"
$srcC.Copy($ws.Range("C309"))
$ws.Range("C309").Value = "HCP"
$ws.Range("D309").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E309"))
$ws.Range("E309").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F309").Value = "compressive yield stress"
$ws.Range("G309").Value = "EXP"
$ws.Range("I309").Value = 298
$ws.Range("J309").Value = 1517000000
$ws.Range("K309").Value = 30000000
$ws.Range("L309").Value = "Pa"
$ws.Range("M309").Value = "T3"
$ws.Range("N309").Value = "10.1016/j.matchar.2024.113730"

# Row 310
$ws.Range("A310").Value = "Al5"
$srcB.Copy($ws.Range("B310"))
$ws.Range("B310").Value = "Al5 Sc25 Ti35 Zr35"
$srcC.Copy($ws.Range("C310"))
$ws.Range("C310").Value = "HCP"
$srcD48.Copy($ws.Range("D310"))
$ws.Range("D310").Value = "AAM"
$ws.Range("F310").Value = "UCS"
$ws.Range("G310").Value = "EXP"
$ws.Range("I310").Value = 298
$ws.Range("J310").Value = 1861000000
$ws.Range("K310").Value = 46000000
$ws.Range("L310").Value = "Pa"
$ws.Range("M310").Value = "T3"
$ws.Range("N310").Value = "10.1016/j.matchar.2024.113730"

# Row 311
$ws.Range("A311").Value = "Al10"
$srcB.Copy($ws.Range("B311"))
$ws.Range("B311").Value = "Al10 Sc20 Ti35 Zr35"
$srcC.Copy($ws.Range("C311"))
$ws.Range("C311").Value = "HCP+BCC"
$srcD48.Copy($ws.Range("D311"))
$ws.Range("D311").Value = "AAM"
$ws.Range("F311").Value = "UCS"
$ws.Range("G311").Value = "EXP"
$ws.Range("I311").Value = 298
$ws.Range("J311").Value = 2040000000
$ws.Range("K311").Value = 25000000
$ws.Range("L311").Value = "Pa"
$ws.Range("M311").Value = "T3"
$ws.Range("N311").Value = "10.1016/j.matchar.2024.113730"

# Row 312
$ws.Range("A312").Value = "Al15"
$srcB.Copy($ws.Range("B312"))
$ws.Range("B312").Value = "Al15 Sc15 Ti35 Zr35"
$srcC.Copy($ws.Range("C312"))
$ws.Range("C312").Value = "BCC"
$srcD48.Copy($ws.Range("D312"))
$ws.Range("D312").Value = "AAM"
$ws.Range("F312").Value = "UCS"
$ws.Range("G312").Value = "EXP"
$ws.Range("I312").Value = 298
$ws.Range("J312").Value = 1798000000
$ws.Range("K312").Value = 146000000
$ws.Range("L312").Value = "Pa"
$ws.Range("M312").Value = "T3"
$ws.Range("N312").Value = "10.1016/j.matchar.2024.113730"

# Row 313
$ws.Range("A313").Value = "Al5"
$srcB.Copy($ws.Range("B313"))
$ws.Range("B313").Value = "Al5 Sc25 Ti35 Zr35"
$srcC.Copy($ws.Range("C313"))
$ws.Range("C313").Value = "HCP+HCP"
$ws.Range("D313").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E313"))
$ws.Range("E313").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F313").Value = "UCS"
$ws.Range("G313").Value = "EXP"
$ws.Range("I313").Value = 298
$ws.Range("J313").Value = 1953000000
$ws.Range("K313").Value = 39000000
$ws.Range("L313").Value = "Pa"
$ws.Range("M313").Value = "T3"
$ws.Range("N313").Value = "10.1016/j.matchar.2024.113730"

# Row 314
$ws.Range("A314").Value = "Al10"
$srcB.Copy($ws.Range("B314"))
$ws.Range("B314").Value = "Al10 Sc20 Ti35 Zr35"
$srcC.Copy($ws.Range("C314"))
$ws.Range("C314").Value = "HCP"
$ws.Range("D314").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E314"))
$ws.Range("E314").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F314").Value = "UCS"
$ws.Range("G314").Value = "EXP"
$ws.Range("I314").Value = 298
$ws.Range("J314").Value = 2584000000
$ws.Range("K314").Value = 44000000
$ws.Range("L314").Value = "Pa"
$ws.Range("M314").Value = "T3"
$ws.Range("N314").Value = "10.1016/j.matchar.2024.113730"

# Row 315
$ws.Range("A315").Value = "Al15"
$srcB.Copy($ws.Range("B315"))
$ws.Range("B315").Value = "Al15 Sc15 Ti35 Zr35"
$srcC.Copy($ws.Range("C315"))
$ws.Range("C315").Value = "BCC+HCP"
$ws.Range("D315").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E315"))
$ws.Range("E315").Value = "annealed in vacuum at 1273K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F315").Value = "UCS"
$ws.Range("G315").Value = "EXP"
$ws.Range("I315").Value = 298
$ws.Range("J315").Value = 1750000000
$ws.Range("K315").Value = 58000000
$ws.Range("L315").Value = "Pa"
$ws.Range("M315").Value = "T3"
$ws.Range("N315").Value = "10.1016/j.matchar.2024.113730"

# Row 316
$ws.Range("A316").Value = "Al5"
$srcB.Copy($ws.Range("B316"))
$ws.Range("B316").Value = "Al5 Sc25 Ti35 Zr35"
$srcC.Copy($ws.Range("C316"))
$ws.Range("C316").Value = "HCP"
$srcD48.Copy($ws.Range("D316"))
$ws.Range("D316").Value = "AAM"
$ws.Range("F316").Value = "compressive ductility"
$ws.Range("G316").Value = "EXP"
$ws.Range("I316").Value = 298
$ws.Range("J316").Value = 4.5
$ws.Range("K316").Value = 0.7
$ws.Range("L316").Value = "%"
$ws.Range("M316").Value = "T3"
$ws.Range("N316").Value = "10.1016/j.matchar.2024.113730"

# Row 317
$ws.Range("A317").Value = "Al10"
$srcB.Copy($ws.Range("B317"))
$ws.Range("B317").Value = "Al10 Sc20 Ti35 Zr35"
$srcC.Copy($ws.Range("C317"))
$ws.Range("C317").Value = "HCP+BCC"
$srcD48.Copy($ws.Range("D317"))
$ws.Range("D317").Value = "AAM"
$ws.Range("F317").Value = "compressive ductility"
$ws.Range("G317").Value = "EXP"
$ws.Range("I317").Value = 298
$ws.Range("J317").Value = 14.9
$ws.Range("K317").Value = 0.7
$ws.Range("L317").Value = "%"
$ws.Range("M317").Value = "T3"
$ws.Range("N317").Value = "10.1016/j.matchar.2024.113730"

# Row 318
$ws.Range("A318").Value = "Al15"
$srcB.Copy($ws.Range("B318"))
$ws.Range("B318").Value = "Al15 Sc15 Ti35 Zr35"
$srcC.Copy($ws.Range("C318"))
$ws.Range("C318").Value = "BCC"
$srcD48.Copy($ws.Range("D318"))
$ws.Range("D318").Value = "AAM"
$ws.Range("F318").Value = "compressive ductility"
$ws.Range("G318").Value = "EXP"
$ws.Range("I318").Value = 298
$ws.Range("J318").Value = 17
$ws.Range("K318").Value = 1.6
$ws.Range("L318").Value = "%"
$ws.Range("M318").Value = "T3"
$ws.Range("N318").Value = "10.1016/j.matchar.2024.113730"

# Row 319
$ws.Range("A319").Value = "Al5"
$srcB.Copy($ws.Range("B319"))
$ws.Range("B319").Value = "Al5 Sc25 Ti35 Zr35"
$srcC.Copy($ws.Range("C319"))
$ws.Range("C319").Value = "HCP+HCP"
$ws.Range("D319").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E319"))
$ws.Range("E319").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F319").Value = "compressive ductility"
$ws.Range("G319").Value = "EXP"
$ws.Range("I319").Value = 298
$ws.Range("J319").Value = 7.8
$ws.Range("K319").Value = 0.6
$ws.Range("L319").Value = "%"
$ws.Range("M319").Value = "T3"
$ws.Range("N319").Value = "10.1016/j.matchar.2024.113730"

# Row 320
$ws.Range("A320").Value = "Al10"
$srcB.Copy($ws.Range("B320"))
$ws.Range("B320").Value = "Al10 Sc20 Ti35 Zr35"
$srcC.Copy($ws.Range("C320"))
$ws.Range("C320").Value = "HCP"
$ws.Range("D320").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E320"))
$ws.Range("E320").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F320").Value = "compressive ductility"
$ws.Range("G320").Value = "EXP"
$ws.Range("I320").Value = 298
$ws.Range("J320").Value = 31.9
$ws.Range("K320").Value = 0.6
$ws.Range("L320").Value = "%"
$ws.Range("M320").Value = "T3"
$ws.Range("N320").Value = "10.1016/j.matchar.2024.113730"

# Row 321
$ws.Range("A321").Value = "Al15"
$srcB.Copy($ws.Range("B321"))
$ws.Range("B321").Value = "Al15 Sc15 Ti35 Zr35"
$srcC.Copy($ws.Range("C321"))
$ws.Range("C321").Value = "BCC+HCP"
$ws.Range("D321").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E321"))
$ws.Range("E321").Value = "annealed in vacuum at 1273K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F321").Value = "compressive ductility"
$ws.Range("G321").Value = "EXP"
$ws.Range("I321").Value = 298
$ws.Range("J321").Value = 2.4
$ws.Range("K321").Value = 0.05
$ws.Range("L321").Value = "%"
$ws.Range("M321").Value = "T3"
$ws.Range("N321").Value = "10.1016/j.matchar.2024.113730"

# Row 322
$ws.Range("A322").Value = "Al5"
$srcB.Copy($ws.Range("B322"))
$ws.Range("B322").Value = "Al5 Sc25 Ti35 Zr35"
$srcC.Copy($ws.Range("C322"))
$ws.Range("C322").Value = "HCP"
$srcD48.Copy($ws.Range("D322"))
$ws.Range("D322").Value = "AAM"
$ws.Range("F322").Value = "density"
$ws.Range("G322").Value = "EXP"
$ws.Range("I322").Value = 298
$ws.Range("J322").Value = 4950
$ws.Range("L322").Value = "kg/m^3"
$ws.Range("M322").Value = "T4"
$ws.Range("N322").Value = "10.1016/j.matchar.2024.113730"

# Row 323
$ws.Range("A323").Value = "Al10"
$srcB.Copy($ws.Range("B323"))
$ws.Range("B323").Value = "Al10 Sc20 Ti35 Zr35"
$srcC.Copy($ws.Range("C323"))
$ws.Range("C323").Value = "HCP+BCC"
$srcD48.Copy($ws.Range("D323"))
$ws.Range("D323").Value = "AAM"
$ws.Range("F323").Value = "density"
$ws.Range("G323").Value = "EXP"
$ws.Range("I323").Value = 298
$ws.Range("J323").Value = 4954
$ws.Range("L323").Value = "kg/m^3"
$ws.Range("M323").Value = "T4"
$ws.Range("N323").Value = "10.1016/j.matchar.2024.113730"

# Row 324
$ws.Range("A324").Value = "Al15"
$srcB.Copy($ws.Range("B324"))
$ws.Range("B324").Value = "Al15 Sc15 Ti35 Zr35"
$srcC.Copy($ws.Range("C324"))
$ws.Range("C324").Value = "BCC"
$srcD48.Copy($ws.Range("D324"))
$ws.Range("D324").Value = "AAM"
$ws.Range("F324").Value = "density"
$ws.Range("G324").Value = "EXP"
$ws.Range("I324").Value = 298
$ws.Range("J324").Value = 4957
$ws.Range("L324").Value = "kg/m^3"
$ws.Range("M324").Value = "T4"
$ws.Range("N324").Value = "10.1016/j.matchar.2024.113730"

# Row 325
$ws.Range("A325").Value = "Al15"
$srcB.Copy($ws.Range("B325"))
$ws.Range("B325").Value = "Al15 Sc15 Ti35 Zr35"
$srcC.Copy($ws.Range("C325"))
$ws.Range("C325").Value = "BCC"
$srcD48.Copy($ws.Range("D325"))
$ws.Range("D325").Value = "AAM"
$ws.Range("F325").Value = "UCS"
$ws.Range("G325").Value = "EXP"
$ws.Range("I325").Value = 873
$ws.Range("J325").Value = 804000000
$ws.Range("L325").Value = "Pa"
$ws.Range("M325").Value = "P5"
$ws.Range("N325").Value = "10.1016/j.matchar.2024.113730"

# Row 326
$ws.Range("A326").Value = "Al5"
$srcB.Copy($ws.Range("B326"))
$ws.Range("B326").Value = "Al5 Sc25 Ti35 Zr35"
$srcC.Copy($ws.Range("C326"))
$ws.Range("C326").Value = "HCP+HCP"
$ws.Range("D326").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E326"))
$ws.Range("E326").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F326").Value = "UCS"
$ws.Range("G326").Value = "EXP"
$ws.Range("I326").Value = 873
$ws.Range("J326").Value = 692000000
$ws.Range("L326").Value = "Pa"
$ws.Range("M326").Value = "P5"
$ws.Range("N326").Value = "10.1016/j.matchar.2024.113730"

# Row 327
$ws.Range("A327").Value = "Al10"
$srcB.Copy($ws.Range("B327"))
$ws.Range("B327").Value = "Al10 Sc20 Ti35 Zr35"
$srcC.Copy($ws.Range("C327"))
$ws.Range("C327").Value = "HCP"
$ws.Range("D327").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E327"))
$ws.Range("E327").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F327").Value = "UCS"
$ws.Range("G327").Value = "EXP"
$ws.Range("I327").Value = 873
$ws.Range("J327").Value = 605000000
$ws.Range("L327").Value = "Pa"
$ws.Range("M327").Value = "P5"
$ws.Range("N327").Value = "10.1016/j.matchar.2024.113730"

# Row 328
$ws.Range("A328").Value = "Al15"
$srcB.Copy($ws.Range("B328"))
$ws.Range("B328").Value = "Al15 Sc15 Ti35 Zr35"
$srcC.Copy($ws.Range("C328"))
$ws.Range("C328").Value = "BCC"
$srcD48.Copy($ws.Range("D328"))
$ws.Range("D328").Value = "AAM"
$ws.Range("F328").Value = "compressive fracture strength"
$ws.Range("G328").Value = "EXP"
$ws.Range("I328").Value = 873
$ws.Range("J328").Value = 689000000
$ws.Range("L328").Value = "Pa"
$ws.Range("M328").Value = "P5"
$ws.Range("N328").Value = "10.1016/j.matchar.2024.113730"

# Row 329
$ws.Range("A329").Value = "Al5"
$srcB.Copy($ws.Range("B329"))
$ws.Range("B329").Value = "Al5 Sc25 Ti35 Zr35"
$srcC.Copy($ws.Range("C329"))
$ws.Range("C329").Value = "HCP+HCP"
$ws.Range("D329").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E329"))
$ws.Range("E329").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F329").Value = "compressive fracture strength"
$ws.Range("G329").Value = "EXP"
$ws.Range("I329").Value = 873
$ws.Range("J329").Value = 576000000
$ws.Range("L329").Value = "Pa"
$ws.Range("M329").Value = "P5"
$ws.Range("N329").Value = "10.1016/j.matchar.2024.113730"

# Row 330
$ws.Range("A330").Value = "Al10"
$srcB.Copy($ws.Range("B330"))
$ws.Range("B330").Value = "Al10 Sc20 Ti35 Zr35"
$srcC.Copy($ws.Range("C330"))
$ws.Range("C330").Value = "HCP"
$ws.Range("D330").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E330"))
$ws.Range("E330").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F330").Value = "compressive fracture strength"
$ws.Range("G330").Value = "EXP"
$ws.Range("I330").Value = 873
$ws.Range("J330").Value = 442000000
$ws.Range("L330").Value = "Pa"
$ws.Range("M330").Value = "P5"
$ws.Range("N330").Value = "10.1016/j.matchar.2024.113730"

# Row 331
$ws.Range("A331").Value = "Al15"
$srcB.Copy($ws.Range("B331"))
$ws.Range("B331").Value = "Al15 Sc15 Ti35 Zr35"
$srcC.Copy($ws.Range("C331"))
$ws.Range("C331").Value = "BCC"
$srcD48.Copy($ws.Range("D331"))
$ws.Range("D331").Value = "AAM"
$ws.Range("F331").Value = "minimum compressive ductility"
$ws.Range("G331").Value = "EXP"
$ws.Range("I331").Value = 873
$ws.Range("J331").Value = 50
$ws.Range("L331").Value = "%"
$ws.Range("M331").Value = "P5"
$ws.Range("N331").Value = "10.1016/j.matchar.2024.113730"

# Row 332
$ws.Range("A332").Value = "Al5"
$srcB.Copy($ws.Range("B332"))
$ws.Range("B332").Value = "Al5 Sc25 Ti35 Zr35"
$srcC.Copy($ws.Range("C332"))
$ws.Range("C332").Value = "HCP+HCP"
$ws.Range("D332").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E332"))
$ws.Range("E332").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F332").Value = "compressive ductility"
$ws.Range("G332").Value = "EXP"
$ws.Range("I332").Value = 873
$ws.Range("J332").Value = 45
$ws.Range("L332").Value = "%"
$ws.Range("M332").Value = "P5"
$ws.Range("N332").Value = "10.1016/j.matchar.2024.113730"

# Row 333
$ws.Range("A333").Value = "Al10"
$srcB.Copy($ws.Range("B333"))
$ws.Range("B333").Value = "Al10 Sc20 Ti35 Zr35"
$srcC.Copy($ws.Range("C333"))
$ws.Range("C333").Value = "HCP"
$ws.Range("D333").Value = "AAM+A+WQ+A"
$srcE.Copy($ws.Range("E333"))
$ws.Range("E333").Value = "annealed in vacuum at 1173K for 60min in quartz then water quenched and annealed again 773K for 30min and cooled down"
$ws.Range("F333").Value = "minimum compressive ductility"
$ws.Range("G333").Value = "EXP"
$ws.Range("I333").Value = 873
$ws.Range("J333").Value = 50
$ws.Range("L333").Value = "%"
$ws.Range("M333").Value = "P5"
$ws.Range("N333").Value = "10.1016/j.matchar.2024.113730"

# Update sheet view / selection to match target state
$ws.Range("A299").Select()
$ws.Range("N337").Select()
